$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 188, shifting all subsequent rows down by one.
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new weekly data record.
$ws.Cells.Item(188, 1).Value  = 4
$ws.Cells.Item(188, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(188, 3).Value  = "Los Lagos"
$ws.Cells.Item(188, 4).Value  = 45258
$ws.Cells.Item(188, 5).Value  = 10
$ws.Cells.Item(188, 6).Value  = 100112052
$ws.Cells.Item(188, 7).Value  = "Albahaca"
$ws.Cells.Item(188, 8).Value  = "Sin especificar"
$ws.Cells.Item(188, 9).Value  = "Primera"
$ws.Cells.Item(188, 10).Value = 90
$ws.Cells.Item(188, 11).Value = 8000
$ws.Cells.Item(188, 12).Value = 8000
$ws.Cells.Item(188, 13).Value = 8000
$ws.Cells.Item(188, 14).Value = "$/docena de matas"
$ws.Cells.Item(188, 15).Value = "Región Metropolitana"
$ws.Cells.Item(188, 16).Value = 1333
$ws.Cells.Item(188, 17).Value = 6
$ws.Cells.Item(188, 18).Value = "Hortaliza"
